$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Medication Class"
$ws.Range("B1").Value = "Classification Code"
$ws.Range("C1").Value = "n"
$ws.Range("D1").Value = "%"

# --- Data rows (re-sorted / refreshed data, V3) ---
$data = @(
    @("ANTIHYPERTENSIVES", 3600000000, 1908, 0.5),
    @("ANALGESICS - NONNARCOTIC", 6400000000, 1719, 0.45),
    @("ANTIHYPERLIPIDEMICS", 3900000000, 1277, 0.33),
    @("MULTIVITAMINS", 7800000000, 1207, 0.32),
    @("DIURETICS", 3700000000, 1052, 0.28),
    @("CALCIUM CHANNEL BLOCKERS", 3400000000, 887, 0.23),
    @("ANTIDIABETICS", 2700000000, 884, 0.23),
    @("MINERALS & ELECTROLYTES", 7900000000, 875, 0.23),
    @("ANALGESICS - ANTI-INFLAMMATORY", 6600000000, 775, 0.2),
    @("BETA BLOCKERS", 3300000000, 711, 0.19),
    @("ULCER DRUGS", 4900000000, 677, 0.18),
    @("VITAMINS", 7700000000, 638, 0.17)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row = $row + 1
}

# --- Formatting ---
# Header row: left align A1, center align B1:D1
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("B1:D1").HorizontalAlignment = -4108

# Body rows: column A left aligned, columns B:C centered, column D percent + centered
$ws.Range("A2:A13").HorizontalAlignment = -4131
$ws.Range("B2:C13").HorizontalAlignment = -4108
$ws.Range("D2:D13").HorizontalAlignment = -4108
$ws.Range("D2:D13").NumberFormat = "0%"
$ws.Range("D2:D13").Style = "Percent"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 33.6328125
$ws.Columns.Item(2).ColumnWidth = 17.7265625
$ws.Columns.Item(3).ColumnWidth = 10.36328125
$ws.Columns.Item(4).ColumnWidth = 10.36328125

# --- Selection ---
$ws.Range("C18").Select()
